$wb = $excel.ActiveWorkbook

# --- Sheet "3 V 0.3": add a new row 9 ---
$ws1 = $wb.Worksheets.Item("3 V 0.3")

$ws1.Cells.Item(9, 1).Value = "20/06/2024 07:45:40"
$ws1.Cells.Item(9, 2).Value = 1
$ws1.Cells.Item(9, 3).Value = "LXCHEM"
$ws1.Cells.Item(9, 4).Value = "Laxmi Organic Industries Ltd"

# E9 must stay a text string ("543277"), not get auto-converted to a number.
$ws1.Cells.Item(9, 5).NumberFormat = "@"
$ws1.Cells.Item(9, 5).Value = "543277"
$ws1.Cells.Item(9, 5).Style = "Normal"

$ws1.Cells.Item(9, 6).Value = 3.2
$ws1.Cells.Item(9, 7).Value = 262.56
$ws1.Cells.Item(9, 8).Value = 4525214

# --- Sheet "DND 3 V 0.3": fix E4 to be numeric, and add new row 5 ---
$ws2 = $wb.Worksheets.Item("DND 3 V 0.3")

$ws2.Cells.Item(4, 5).Value = 532832

$ws2.Cells.Item(5, 1).Value = "20/06/2024 07:45:40"
$ws2.Cells.Item(5, 2).Value = 1
$ws2.Cells.Item(5, 3).Value = "IBREALEST"
$ws2.Cells.Item(5, 4).Value = "Indiabulls Real Estate Limited"

# E5 must stay a text string ("532832"), not get auto-converted to a number.
$ws2.Cells.Item(5, 5).NumberFormat = "@"
$ws2.Cells.Item(5, 5).Value = "532832"
$ws2.Cells.Item(5, 5).Style = "Normal"

$ws2.Cells.Item(5, 6).Value = 14.08
$ws2.Cells.Item(5, 7).Value = 156.2
$ws2.Cells.Item(5, 8).Value = 75286146
